$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new "Media" column: shift old "Fixed Income" header from M1 to N1,
# put "Media" in M1, and copy the header style (bold/border/centered) to N1.
$ws.Cells.Item(1, 13).Value = "Media"
$ws.Cells.Item(1, 14).Value = "Fixed Income"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# Row 2: updated trade returns / closing prices
$ws.Cells.Item(2, 3).Value = -0.1489810950416438
$ws.Cells.Item(2, 4).Value = -0.2043474562481407
$ws.Cells.Item(2, 5).Value = -0.2547938095604285
$ws.Cells.Item(2, 6).Value = -0.0283171465889092
$ws.Cells.Item(2, 7).Value = -1.561687929733194
$ws.Cells.Item(2, 8).Value = 0.2521240374853275
$ws.Cells.Item(2, 9).Value = -0.07928124633826795
$ws.Cells.Item(2, 10).Value = 0.02041826591931516
$ws.Cells.Item(2, 11).Value = 0.02378011117973977
$ws.Cells.Item(2, 12).Value = -0.2098342518997854
$ws.Cells.Item(2, 13).Value = -0.5453174658749476
$ws.Cells.Item(2, 14).Value = 0.260637167218768

# Row 3: updated trade returns / closing prices
$ws.Cells.Item(3, 3).Value = -0.4062303152217909
$ws.Cells.Item(3, 4).Value = -0.03086517683768375
$ws.Cells.Item(3, 5).Value = -0.2300319671529607
$ws.Cells.Item(3, 6).Value = -0.03032195891464575
$ws.Cells.Item(3, 7).Value = -1.186312698696654
$ws.Cells.Item(3, 8).Value = 0.3629622412102476
$ws.Cells.Item(3, 9).Value = -0.5040763229969668
$ws.Cells.Item(3, 10).Value = 0.6147778012059749
$ws.Cells.Item(3, 11).Value = -0.2386391904336048
$ws.Cells.Item(3, 12).Value = -2.109982518704936
$ws.Cells.Item(3, 13).Value = 0.3768739117126277
$ws.Cells.Item(3, 14).Value = 0.8037272180612992

# Row 4: updated trade returns / closing prices
$ws.Cells.Item(4, 3).Value = -0.7169975486779477
$ws.Cells.Item(4, 4).Value = -0.2946698798083034
$ws.Cells.Item(4, 5).Value = -0.06149789202139089
$ws.Cells.Item(4, 6).Value = 0.01091433704255712
$ws.Cells.Item(4, 7).Value = 0.306022683792049
$ws.Cells.Item(4, 8).Value = 0.1005380592722608
$ws.Cells.Item(4, 9).Value = -3.047512040033329
$ws.Cells.Item(4, 10).Value = -0.1497267621519438
$ws.Cells.Item(4, 11).Value = -0.02892383927239406
$ws.Cells.Item(4, 12).Value = -2.592559864161619
$ws.Cells.Item(4, 13).Value = 0.2672110784188562
$ws.Cells.Item(4, 14).Value = -0.7978651073472509

# Row 5: updated trade returns / closing prices
$ws.Cells.Item(5, 3).Value = -0.5414142807543174
$ws.Cells.Item(5, 4).Value = -0.007461755146021876
$ws.Cells.Item(5, 5).Value = -0.08557828813120216
$ws.Cells.Item(5, 6).Value = 0.274611466321977
$ws.Cells.Item(5, 7).Value = -0.7103281350806441
$ws.Cells.Item(5, 8).Value = -0.1362701662175447
$ws.Cells.Item(5, 9).Value = -2.945418105643589
$ws.Cells.Item(5, 10).Value = -0.08326532583783304
$ws.Cells.Item(5, 11).Value = -0.2787265467134047
$ws.Cells.Item(5, 12).Value = -4.871307790857954
$ws.Cells.Item(5, 13).Value = 0.7994323085619569
$ws.Cells.Item(5, 14).Value = 0.604667271479058

# Row 6: updated trade returns / closing prices
$ws.Cells.Item(6, 3).Value = -1.273814175376051
$ws.Cells.Item(6, 4).Value = 0.8625605506452008
$ws.Cells.Item(6, 5).Value = -0.3516343311900239
$ws.Cells.Item(6, 6).Value = 0.8295064829069797
$ws.Cells.Item(6, 7).Value = 0.2646541652774737
$ws.Cells.Item(6, 8).Value = -2.295526176009214
$ws.Cells.Item(6, 9).Value = -5.359240583353815
$ws.Cells.Item(6, 10).Value = -0.01152431248113572
$ws.Cells.Item(6, 11).Value = -1.349267970439772
$ws.Cells.Item(6, 12).Value = -10.56842778552956
$ws.Cells.Item(6, 13).Value = 3.943185013355406
$ws.Cells.Item(6, 14).Value = -1.726398582041249

# Row 7: updated trade returns / closing prices
$ws.Cells.Item(7, 3).Value = -0.2142419044454809
$ws.Cells.Item(7, 4).Value = 0.1503728058796788
$ws.Cells.Item(7, 5).Value = -0.05601791128669487
$ws.Cells.Item(7, 6).Value = 0.1404504118786588
$ws.Cells.Item(7, 7).Value = 0.04807982169657055
$ws.Cells.Item(7, 8).Value = -0.3776253217542895
$ws.Cells.Item(7, 9).Value = -1.049080538168349
$ws.Cells.Item(7, 10).Value = 0.007488712061596708
$ws.Cells.Item(7, 11).Value = -0.2382796914939904
$ws.Cells.Item(7, 12).Value = -1.680918564179646
$ws.Cells.Item(7, 13).Value = 0.6853788724589751
$ws.Cells.Item(7, 14).Value = -0.3174670141187252

# Row 8: updated trade returns / closing prices
$ws.Cells.Item(8, 3).Value = 0.01773549883990645
$ws.Cells.Item(8, 4).Value = 0.03344697747655328
$ws.Cells.Item(8, 5).Value = -0.2902820355610058
$ws.Cells.Item(8, 6).Value = -0.01116174402250351
$ws.Cells.Item(8, 7).Value = -0.01692836908134405
$ws.Cells.Item(8, 8).Value = 0.02704656276890775
$ws.Cells.Item(8, 9).Value = -0.000989563200618839
$ws.Cells.Item(8, 10).Value = -0.04871881893772213
$ws.Cells.Item(8, 11).Value = 0.02378011117973977
$ws.Cells.Item(8, 12).Value = 0.1092988668555154
$ws.Cells.Item(8, 13).Value = -0.007527601204413203
$ws.Cells.Item(8, 14).Value = 0

# Row 9: updated trade returns / closing prices
$ws.Cells.Item(9, 3).Value = 0.1955532435975876
$ws.Cells.Item(9, 4).Value = 0.2316634406487302
$ws.Cells.Item(9, 5).Value = -0.4144248808142237
$ws.Cells.Item(9, 6).Value = -0.02733816106486643
$ws.Cells.Item(9, 7).Value = 0.01173841904374701
$ws.Cells.Item(9, 8).Value = -0.05855074679567809
$ws.Cells.Item(9, 9).Value = -0.08736632778662594
$ws.Cells.Item(9, 10).Value = 0.002565935620940141
$ws.Cells.Item(9, 11).Value = -0.2386391904336048
$ws.Cells.Item(9, 12).Value = -1.274290772393061
$ws.Cells.Item(9, 13).Value = 0.6701319963834163
$ws.Cells.Item(9, 14).Value = 0

# Row 10: updated trade returns / closing prices
$ws.Cells.Item(10, 3).Value = 0.09140145510292352
$ws.Cells.Item(10, 4).Value = 0.2905202403063736
$ws.Cells.Item(10, 5).Value = -1.317140112042914
$ws.Cells.Item(10, 6).Value = -0.09078809605400978
$ws.Cells.Item(10, 7).Value = -0.2017371767166316
$ws.Cells.Item(10, 8).Value = -0.04500878434542416
$ws.Cells.Item(10, 9).Value = 0.03521847168395729
$ws.Cells.Item(10, 10).Value = -0.150905054872121
$ws.Cells.Item(10, 11).Value = -0.02892383927239406
$ws.Cells.Item(10, 12).Value = -2.615659896719902
$ws.Cells.Item(10, 13).Value = 0.9646051872465706
$ws.Cells.Item(10, 14).Value = 0

# Row 11: updated trade returns / closing prices
$ws.Cells.Item(11, 3).Value = 0.1909877564979118
$ws.Cells.Item(11, 4).Value = 0.3282389310722448
$ws.Cells.Item(11, 5).Value = -1.224729385953708
$ws.Cells.Item(11, 6).Value = -0.06429085245209593
$ws.Cells.Item(11, 7).Value = -0.1102521148231073
$ws.Cells.Item(11, 8).Value = -0.186738143213691
$ws.Cells.Item(11, 9).Value = -0.0005636085943564948
$ws.Cells.Item(11, 10).Value = -0.06378229257892326
$ws.Cells.Item(11, 11).Value = -0.2787265467134047
$ws.Cells.Item(11, 12).Value = -3.915233991133615
$ws.Cells.Item(11, 13).Value = 1.349852520276009
$ws.Cells.Item(11, 14).Value = 0

# Row 12: updated trade returns / closing prices
$ws.Cells.Item(12, 3).Value = 0.8386172526222808
$ws.Cells.Item(12, 4).Value = 1.027109951478823
$ws.Cells.Item(12, 5).Value = -0.9211664984126267
$ws.Cells.Item(12, 6).Value = 0.3607554133549377
$ws.Cells.Item(12, 7).Value = 0.2875776170783487
$ws.Cells.Item(12, 8).Value = -0.5264940492721689
$ws.Cells.Item(12, 9).Value = -0.4025158197733269
$ws.Cells.Item(12, 10).Value = 0.754839523181165
$ws.Cells.Item(12, 11).Value = -1.349267970439772
$ws.Cells.Item(12, 12).Value = -9.061371058999205
$ws.Cells.Item(12, 13).Value = 1.624013827346932
$ws.Cells.Item(12, 14).Value = 0

# Row 13: updated trade returns / closing prices
$ws.Cells.Item(13, 3).Value = 0.1493464753321933
$ws.Cells.Item(13, 4).Value = 0.1742557189694822
$ws.Cells.Item(13, 5).Value = -0.2479797273564791
$ws.Cells.Item(13, 6).Value = 0.06583425231915273
$ws.Cells.Item(13, 7).Value = 0.05178587098626595
$ws.Cells.Item(13, 8).Value = -0.0922504335250682
$ws.Cells.Item(13, 9).Value = -0.07035826605024999
$ws.Cells.Item(13, 10).Value = 0.1334791526336825
$ws.Cells.Item(13, 11).Value = -0.2382796914939904
$ws.Cells.Item(13, 12).Value = -1.47405577885728
$ws.Cells.Item(13, 13).Value = 0.3066904066269789
$ws.Cells.Item(13, 14).Value = 0

# Row 14: updated trade returns / closing prices
$ws.Cells.Item(14, 3).Value = -0.1667165938815502
$ws.Cells.Item(14, 4).Value = -0.237794433724694
$ws.Cells.Item(14, 5).Value = 0.03548822600057738
$ws.Cells.Item(14, 6).Value = -0.01715540256640569
$ws.Cells.Item(14, 7).Value = -1.54475956065185
$ws.Cells.Item(14, 8).Value = 0.2250774747164198
$ws.Cells.Item(14, 9).Value = -0.07829168313764912
$ws.Cells.Item(14, 10).Value = 0.06913708485703729
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = -0.3191331187553008
$ws.Cells.Item(14, 13).Value = -0.5377898646705344
$ws.Cells.Item(14, 14).Value = 0.260637167218768

# Row 15: updated trade returns / closing prices
$ws.Cells.Item(15, 3).Value = -0.6017835588193785
$ws.Cells.Item(15, 4).Value = -0.2625286174864139
$ws.Cells.Item(15, 5).Value = 0.184392913661263
$ws.Cells.Item(15, 6).Value = -0.00298379784977932
$ws.Cells.Item(15, 7).Value = -1.198051117740401
$ws.Cells.Item(15, 8).Value = 0.4215129880059256
$ws.Cells.Item(15, 9).Value = -0.4167099952103409
$ws.Cells.Item(15, 10).Value = 0.6122118655850348
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = -0.8356917463118756
$ws.Cells.Item(15, 13).Value = -0.2932580846707885
$ws.Cells.Item(15, 14).Value = 0.8037272180612992

# Row 16: updated trade returns / closing prices
$ws.Cells.Item(16, 3).Value = -0.8083990037808713
$ws.Cells.Item(16, 4).Value = -0.585190120114677
$ws.Cells.Item(16, 5).Value = 1.255642220021523
$ws.Cells.Item(16, 6).Value = 0.1017024330965669
$ws.Cells.Item(16, 7).Value = 0.5077598605086806
$ws.Cells.Item(16, 8).Value = 0.145546843617685
$ws.Cells.Item(16, 9).Value = -3.082730511717286
$ws.Cells.Item(16, 10).Value = 0.001178292720177221
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0.02310003255828352
$ws.Cells.Item(16, 13).Value = -0.6973941088277145
$ws.Cells.Item(16, 14).Value = -0.7978651073472509

# Row 17: updated trade returns / closing prices
$ws.Cells.Item(17, 3).Value = -0.7324020372522293
$ws.Cells.Item(17, 4).Value = -0.3357006862182667
$ws.Cells.Item(17, 5).Value = 1.139151097822506
$ws.Cells.Item(17, 6).Value = 0.3389023187740729
$ws.Cells.Item(17, 7).Value = -0.6000760202575367
$ws.Cells.Item(17, 8).Value = 0.05046797699614627
$ws.Cells.Item(17, 9).Value = -2.944854497049232
$ws.Cells.Item(17, 10).Value = -0.01948303325890978
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = -0.9560737997243387
$ws.Cells.Item(17, 13).Value = -0.5504202117140521
$ws.Cells.Item(17, 14).Value = 0.604667271479058

# Row 18: updated trade returns / closing prices
$ws.Cells.Item(18, 3).Value = -2.112431427998332
$ws.Cells.Item(18, 4).Value = -0.1645494008336221
$ws.Cells.Item(18, 5).Value = 0.5695321672226028
$ws.Cells.Item(18, 6).Value = 0.468751069552042
$ws.Cells.Item(18, 7).Value = -0.02292345180087502
$ws.Cells.Item(18, 8).Value = -1.769032126737045
$ws.Cells.Item(18, 9).Value = -4.956724763580488
$ws.Cells.Item(18, 10).Value = -0.7663638356623007
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = -1.507056726530351
$ws.Cells.Item(18, 13).Value = 2.319171186008474
$ws.Cells.Item(18, 14).Value = -1.726398582041249

# Row 19: updated trade returns / closing prices
$ws.Cells.Item(19, 3).Value = -0.3635883797776742
$ws.Cells.Item(19, 4).Value = -0.02388291308980334
$ws.Cells.Item(19, 5).Value = 0.1919618160697842
$ws.Cells.Item(19, 6).Value = 0.07461615955950603
$ws.Cells.Item(19, 7).Value = -0.003706049289695403
$ws.Cells.Item(19, 8).Value = -0.2853748882292213
$ws.Cells.Item(19, 9).Value = -0.9787222721180988
$ws.Cells.Item(19, 10).Value = 0.007488712061596708
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = -0.2068627853223666
$ws.Cells.Item(19, 13).Value = 0.3786884658319962
$ws.Cells.Item(19, 14).Value = -0.3174670141187252
